$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing column-A pattern (rows 1-195 all hold 0) down
# through row 200, matching rows 196-200 added in the diff.
for ($r = 196; $r -le 200; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
}
